$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at D (old "Avg" column D shifts right to E,
#        carrying its data/style along). ------------------------------------
$ws.Range("D1").EntireColumn.Insert()

# --- 2. Seed the new shared strings in the same order the original commit
#        added them ("Runden:" then "BlocksFitting") by writing the G1
#        label first. ---------------------------------------------------------
$ws.Range("G1").Value = "Runden:"
$ws.Range("C1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# --- 3. New column D: "BlocksFitting" header + values, styled like the
#        neighbouring "singleBlockObs" column (C). ---------------------------
$ws.Range("D1").Value = "BlocksFitting"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("D2").Value = 0
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("D3").Value = 0
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("D4").Value = 0
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column D width, matching column C's.
$ws.Columns.Item(4).ColumnWidth = 13.592447916666666

# --- 4. New data values for rows 3 & 4 (previously blank) and the shifted
#        "Avg" column (now E). ------------------------------------------------
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 0

$ws.Range("E2").Value = 655
$ws.Range("E3").Value = 600

# --- 5. Move the color-scale conditional formatting from D2:D9 to E2:E9
#        (it stayed anchored on D after the column insert). ------------------
$ws.Range("D2:D9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E9"))

# --- 6. "500" value in column H, next to "Runden:". --------------------------
$ws.Range("H1").Value = 500
$ws.Range("A2").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 7. Selection ends on G4, matching the saved view. -----------------------
$ws.Range("G4").Select()
